$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.670.53"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "2.467.44"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0853"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").Value = "2.847.71"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").Value = "2.457.90"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "41.622.84"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  -2.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.72%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0763"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("E39").Value = "  +1.16%  "

$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.82%  "

$ws.Range("D43").Value = "1.998.50"
$ws.Range("E43").Value = "  +1.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.21%  "

$ws.Range("D48").Value = "2.727.64"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.52%  "
